$wb = $excel.ActiveWorkbook

$wsGlobal = $wb.Worksheets.Item("Global")
$originalActive = $wb.ActiveSheet

# Switch the launched browser from Firefox to Chrome
$wsGlobal.Range("A2").Value = "CHROME"

# Move the saved active cell selection on the Global sheet down to A3
$wsGlobal.Range("A3").Select()

# Restore the originally active sheet/tab so it is not changed by the selection above
$originalActive.Activate()
